$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B (target OOXML width 29.7109375; closest representable value)
$ws.Columns.Item(2).ColumnWidth = 28.75

# New row 7 data
$ws.Range("A7").Value = 43058
$ws.Range("B7").Value = "Alles Da.jpg; Angekommen.jpeg"
$ws.Range("C7").Value = "Die Teile sind alle Da! Viel früher als erwartet. Auf die Chinesen ist verlass! Auf zum Löten (Raspi + Kamera fehlen)"

# Reuse the same date formatting as the other date cells in column A
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update selection as in the diff (selection moved to C8)
$ws.Range("C8").Select()
